$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 202, shifting existing rows 202:211 down to 203:212
$ws.Rows.Item(202).Insert()

# Populate the newly inserted row 202 with the new data record
$ws.Cells.Item(202, 1).Value = 5
$ws.Cells.Item(202, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(202, 3).Value = "Maule"
$ws.Cells.Item(202, 4).Value = 44939
$ws.Cells.Item(202, 5).Value = 7
$ws.Cells.Item(202, 6).Value = 100112031
$ws.Cells.Item(202, 7).Value = "Poroto verde"
$ws.Cells.Item(202, 8).Value = "Sin especificar"
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 300
$ws.Cells.Item(202, 11).Value = 28000
$ws.Cells.Item(202, 12).Value = 28000
$ws.Cells.Item(202, 13).Value = 28000
$ws.Cells.Item(202, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(202, 15).Value = "Región del Maule"
$ws.Cells.Item(202, 16).Value = 1120
$ws.Cells.Item(202, 17).Value = 25
$ws.Cells.Item(202, 18).Value = "Hortaliza"
